$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 12: CSK vs RR - fill in scores for each player (row 24)
$ws.Range("E24").Value = 80
$ws.Range("H24").Value = 0
$ws.Range("K24").Value = 60
$ws.Range("N24").Value = 100
$ws.Range("Q24").Value = 80
$ws.Range("T24").Value = 40
$ws.Range("W24").Value = 30
$ws.Range("Z24").Value = 50
$ws.Range("AC24").Value = 20

# Tied scores (E24 == Q24 == 80) -> points are manually split/averaged, overwriting formula
$ws.Range("D24").Value = 12.5
$ws.Range("P24").Value = 12.5
